$d = $word.ActiveDocument

# Direct cell addressing for duplicated "old" text values (Find/Replace All would
# incorrectly overwrite both occurrences with the same replacement).
$t = $d.Tables.Item(1)
$t.Cell(5,5).Range.Text = "37+9=46"  # was "13-7=6"
$t.Cell(9,4).Range.Text = "49+22=71"  # was "13-7=6"
$t.Cell(13,3).Range.Text = "40-33=7"  # was "69+5=74"
$t.Cell(16,3).Range.Text = "81-2=79"  # was "69+5=74"

# Remaining (unique) replacements via Find/Replace across the whole document.
$d.Content.Find.Execute("2025-10-15 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-10-16 Thursday", 2) | Out-Null
$d.Content.Find.Execute("49+34=83", $true, $false, $false, $false, $false, $true, 1, $false, "98-19=79", 2) | Out-Null
$d.Content.Find.Execute("27+47=74", $true, $false, $false, $false, $false, $true, 1, $false, "56+26=82", 2) | Out-Null
$d.Content.Find.Execute("39+39=78", $true, $false, $false, $false, $false, $true, 1, $false, "23+8=31", 2) | Out-Null
$d.Content.Find.Execute("85-47=38", $true, $false, $false, $false, $false, $true, 1, $false, "27+27=54", 2) | Out-Null
$d.Content.Find.Execute("40-23=17", $true, $false, $false, $false, $false, $true, 1, $false, "91-34=57", 2) | Out-Null
$d.Content.Find.Execute("16+8=24", $true, $false, $false, $false, $false, $true, 1, $false, "90-62=28", 2) | Out-Null
$d.Content.Find.Execute("17+49=66", $true, $false, $false, $false, $false, $true, 1, $false, "22-19=3", 2) | Out-Null
$d.Content.Find.Execute("51-7=44", $true, $false, $false, $false, $false, $true, 1, $false, "9+2=11", 2) | Out-Null
$d.Content.Find.Execute("75+18=93", $true, $false, $false, $false, $false, $true, 1, $false, "29+48=77", 2) | Out-Null
$d.Content.Find.Execute("48+13=61", $true, $false, $false, $false, $false, $true, 1, $false, "90-79=11", 2) | Out-Null
$d.Content.Find.Execute("28+17=45", $true, $false, $false, $false, $false, $true, 1, $false, "81-44=37", 2) | Out-Null
$d.Content.Find.Execute("17+77=94", $true, $false, $false, $false, $false, $true, 1, $false, "63-25=38", 2) | Out-Null
$d.Content.Find.Execute("60-1=59", $true, $false, $false, $false, $false, $true, 1, $false, "13+59=72", 2) | Out-Null
$d.Content.Find.Execute("58+23=81", $true, $false, $false, $false, $false, $true, 1, $false, "41-17=24", 2) | Out-Null
$d.Content.Find.Execute("96-48=48", $true, $false, $false, $false, $false, $true, 1, $false, "7+49=56", 2) | Out-Null
$d.Content.Find.Execute("53-14=39", $true, $false, $false, $false, $false, $true, 1, $false, "80-16=64", 2) | Out-Null
$d.Content.Find.Execute("71-47=24", $true, $false, $false, $false, $false, $true, 1, $false, "74-37=37", 2) | Out-Null
$d.Content.Find.Execute("57+39=96", $true, $false, $false, $false, $false, $true, 1, $false, "34+27=61", 2) | Out-Null
$d.Content.Find.Execute("8+29=37", $true, $false, $false, $false, $false, $true, 1, $false, "34-17=17", 2) | Out-Null
$d.Content.Find.Execute("43-34=9", $true, $false, $false, $false, $false, $true, 1, $false, "28+68=96", 2) | Out-Null
$d.Content.Find.Execute("80-51=29", $true, $false, $false, $false, $false, $true, 1, $false, "3+89=92", 2) | Out-Null
$d.Content.Find.Execute("70-63=7", $true, $false, $false, $false, $false, $true, 1, $false, "23+38=61", 2) | Out-Null
$d.Content.Find.Execute("46+46=92", $true, $false, $false, $false, $false, $true, 1, $false, "69+7=76", 2) | Out-Null
$d.Content.Find.Execute("28+34=62", $true, $false, $false, $false, $false, $true, 1, $false, "81-25=56", 2) | Out-Null
$d.Content.Find.Execute("83-15=68", $true, $false, $false, $false, $false, $true, 1, $false, "42-35=7", 2) | Out-Null
$d.Content.Find.Execute("90-73=17", $true, $false, $false, $false, $false, $true, 1, $false, "70-24=46", 2) | Out-Null
$d.Content.Find.Execute("91-6=85", $true, $false, $false, $false, $false, $true, 1, $false, "71-59=12", 2) | Out-Null
$d.Content.Find.Execute("32-25=7", $true, $false, $false, $false, $false, $true, 1, $false, "39+57=96", 2) | Out-Null
$d.Content.Find.Execute("29+9=38", $true, $false, $false, $false, $false, $true, 1, $false, "66+27=93", 2) | Out-Null
$d.Content.Find.Execute("56+15=71", $true, $false, $false, $false, $false, $true, 1, $false, "71-33=38", 2) | Out-Null
$d.Content.Find.Execute("80-25=55", $true, $false, $false, $false, $false, $true, 1, $false, "9+67=76", 2) | Out-Null
$d.Content.Find.Execute("73+18=91", $true, $false, $false, $false, $false, $true, 1, $false, "18+63=81", 2) | Out-Null
$d.Content.Find.Execute("20-5=15", $true, $false, $false, $false, $false, $true, 1, $false, "59+2=61", 2) | Out-Null
$d.Content.Find.Execute("9+66=75", $true, $false, $false, $false, $false, $true, 1, $false, "26+55=81", 2) | Out-Null
$d.Content.Find.Execute("8+75=83", $true, $false, $false, $false, $false, $true, 1, $false, "54+29=83", 2) | Out-Null
$d.Content.Find.Execute("64-56=8", $true, $false, $false, $false, $false, $true, 1, $false, "53-45=8", 2) | Out-Null
$d.Content.Find.Execute("76-37=39", $true, $false, $false, $false, $false, $true, 1, $false, "27-8=19", 2) | Out-Null
$d.Content.Find.Execute("16+16=32", $true, $false, $false, $false, $false, $true, 1, $false, "41-16=25", 2) | Out-Null
$d.Content.Find.Execute("9+4=13", $true, $false, $false, $false, $false, $true, 1, $false, "18+28=46", 2) | Out-Null
$d.Content.Find.Execute("91-79=12", $true, $false, $false, $false, $false, $true, 1, $false, "9+39=48", 2) | Out-Null
$d.Content.Find.Execute("17+7=24", $true, $false, $false, $false, $false, $true, 1, $false, "19+77=96", 2) | Out-Null
$d.Content.Find.Execute("83-68=15", $true, $false, $false, $false, $false, $true, 1, $false, "94-58=36", 2) | Out-Null
$d.Content.Find.Execute("82-25=57", $true, $false, $false, $false, $false, $true, 1, $false, "27+66=93", 2) | Out-Null
$d.Content.Find.Execute("92-83=9", $true, $false, $false, $false, $false, $true, 1, $false, "80-29=51", 2) | Out-Null
$d.Content.Find.Execute("81-19=62", $true, $false, $false, $false, $false, $true, 1, $false, "63+28=91", 2) | Out-Null
$d.Content.Find.Execute("33-16=17", $true, $false, $false, $false, $false, $true, 1, $false, "17+38=55", 2) | Out-Null
$d.Content.Find.Execute("62-19=43", $true, $false, $false, $false, $false, $true, 1, $false, "16+47=63", 2) | Out-Null
$d.Content.Find.Execute("98-49=49", $true, $false, $false, $false, $false, $true, 1, $false, "61-2=59", 2) | Out-Null
$d.Content.Find.Execute("73-27=46", $true, $false, $false, $false, $false, $true, 1, $false, "78+9=87", 2) | Out-Null
$d.Content.Find.Execute("4+89=93", $true, $false, $false, $false, $false, $true, 1, $false, "86-69=17", 2) | Out-Null
$d.Content.Find.Execute("60-48=12", $true, $false, $false, $false, $false, $true, 1, $false, "58-29=29", 2) | Out-Null
$d.Content.Find.Execute("92-45=47", $true, $false, $false, $false, $false, $true, 1, $false, "18+45=63", 2) | Out-Null
$d.Content.Find.Execute("43-38=5", $true, $false, $false, $false, $false, $true, 1, $false, "71-22=49", 2) | Out-Null
$d.Content.Find.Execute("54-47=7", $true, $false, $false, $false, $false, $true, 1, $false, "77+8=85", 2) | Out-Null
$d.Content.Find.Execute("52-9=43", $true, $false, $false, $false, $false, $true, 1, $false, "94-56=38", 2) | Out-Null
$d.Content.Find.Execute("98-79=19", $true, $false, $false, $false, $false, $true, 1, $false, "80-77=3", 2) | Out-Null
$d.Content.Find.Execute("40-16=24", $true, $false, $false, $false, $false, $true, 1, $false, "14+38=52", 2) | Out-Null
$d.Content.Find.Execute("35-26=9", $true, $false, $false, $false, $false, $true, 1, $false, "66+29=95", 2) | Out-Null
$d.Content.Find.Execute("59+29=88", $true, $false, $false, $false, $false, $true, 1, $false, "18+48=66", 2) | Out-Null
$d.Content.Find.Execute("68+17=85", $true, $false, $false, $false, $false, $true, 1, $false, "82-66=16", 2) | Out-Null
$d.Content.Find.Execute("65-17=48", $true, $false, $false, $false, $false, $true, 1, $false, "46+29=75", 2) | Out-Null
$d.Content.Find.Execute("38+43=81", $true, $false, $false, $false, $false, $true, 1, $false, "29+57=86", 2) | Out-Null
$d.Content.Find.Execute("15+47=62", $true, $false, $false, $false, $false, $true, 1, $false, "86-37=49", 2) | Out-Null
$d.Content.Find.Execute("15+77=92", $true, $false, $false, $false, $false, $true, 1, $false, "51-28=23", 2) | Out-Null
$d.Content.Find.Execute("76+5=81", $true, $false, $false, $false, $false, $true, 1, $false, "39+28=67", 2) | Out-Null
$d.Content.Find.Execute("70-65=5", $true, $false, $false, $false, $false, $true, 1, $false, "56+28=84", 2) | Out-Null
$d.Content.Find.Execute("95-58=37", $true, $false, $false, $false, $false, $true, 1, $false, "22-19=3", 2) | Out-Null
$d.Content.Find.Execute("50-48=2", $true, $false, $false, $false, $false, $true, 1, $false, "37+35=72", 2) | Out-Null
$d.Content.Find.Execute("40-2=38", $true, $false, $false, $false, $false, $true, 1, $false, "67+4=71", 2) | Out-Null
$d.Content.Find.Execute("80-7=73", $true, $false, $false, $false, $false, $true, 1, $false, "29+46=75", 2) | Out-Null
$d.Content.Find.Execute("31-17=14", $true, $false, $false, $false, $false, $true, 1, $false, "90-36=54", 2) | Out-Null
$d.Content.Find.Execute("47+45=92", $true, $false, $false, $false, $false, $true, 1, $false, "48-19=29", 2) | Out-Null
$d.Content.Find.Execute("42-5=37", $true, $false, $false, $false, $false, $true, 1, $false, "57+25=82", 2) | Out-Null
$d.Content.Find.Execute("96-29=67", $true, $false, $false, $false, $false, $true, 1, $false, "98-59=39", 2) | Out-Null
$d.Content.Find.Execute("25+16=41", $true, $false, $false, $false, $false, $true, 1, $false, "27+55=82", 2) | Out-Null
$d.Content.Find.Execute("76-49=27", $true, $false, $false, $false, $false, $true, 1, $false, "64-47=17", 2) | Out-Null
$d.Content.Find.Execute("5+38=43", $true, $false, $false, $false, $false, $true, 1, $false, "49+37=86", 2) | Out-Null
$d.Content.Find.Execute("91-58=33", $true, $false, $false, $false, $false, $true, 1, $false, "58+4=62", 2) | Out-Null
$d.Content.Find.Execute("57+28=85", $true, $false, $false, $false, $false, $true, 1, $false, "41-15=26", 2) | Out-Null
$d.Content.Find.Execute("83-5=78", $true, $false, $false, $false, $false, $true, 1, $false, "70-27=43", 2) | Out-Null
$d.Content.Find.Execute("58+39=97", $true, $false, $false, $false, $false, $true, 1, $false, "48+16=64", 2) | Out-Null
$d.Content.Find.Execute("75-6=69", $true, $false, $false, $false, $false, $true, 1, $false, "26+5=31", 2) | Out-Null
$d.Content.Find.Execute("70-6=64", $true, $false, $false, $false, $false, $true, 1, $false, "68-59=9", 2) | Out-Null
$d.Content.Find.Execute("18+55=73", $true, $false, $false, $false, $false, $true, 1, $false, "17+47=64", 2) | Out-Null
$d.Content.Find.Execute("71-66=5", $true, $false, $false, $false, $false, $true, 1, $false, "57-39=18", 2) | Out-Null
$d.Content.Find.Execute("30-9=21", $true, $false, $false, $false, $false, $true, 1, $false, "28+43=71", 2) | Out-Null
$d.Content.Find.Execute("34-26=8", $true, $false, $false, $false, $false, $true, 1, $false, "47-38=9", 2) | Out-Null
$d.Content.Find.Execute("84-65=19", $true, $false, $false, $false, $false, $true, 1, $false, "91-38=53", 2) | Out-Null
$d.Content.Find.Execute("57+5=62", $true, $false, $false, $false, $false, $true, 1, $false, "73-58=15", 2) | Out-Null
$d.Content.Find.Execute("23+49=72", $true, $false, $false, $false, $false, $true, 1, $false, "58+18=76", 2) | Out-Null
$d.Content.Find.Execute("90-23=67", $true, $false, $false, $false, $false, $true, 1, $false, "26-18=8", 2) | Out-Null
$d.Content.Find.Execute("27+37=64", $true, $false, $false, $false, $false, $true, 1, $false, "70-2=68", 2) | Out-Null
$d.Content.Find.Execute("69+14=83", $true, $false, $false, $false, $false, $true, 1, $false, "59+38=97", 2) | Out-Null
$d.Content.Find.Execute("65-47=18", $true, $false, $false, $false, $false, $true, 1, $false, "64+27=91", 2) | Out-Null
$d.Content.Find.Execute("42-8=34", $true, $false, $false, $false, $false, $true, 1, $false, "63-37=26", 2) | Out-Null

Write-Output "done"
